$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 21:48:40"
$ws.Range("E3").Value = "2026-02-07 21:48:43"
$ws.Range("H3").Value = "'86%"
$ws.Range("E4").Value = "2026-02-07 21:48:46"
$ws.Range("N4").Value = "8.9 °C 21:17 TU"
$ws.Range("O4").Value = "11.9 °C"
$ws.Range("E5").Value = "2026-02-07 21:48:48"
$ws.Range("E6").Value = "2026-02-07 21:48:51"
$ws.Range("E7").Value = "2026-02-07 21:48:54"
$ws.Range("O7").Value = "13.1 °C"
$ws.Range("E8").Value = "2026-02-07 21:48:56"
$ws.Range("E9").Value = "2026-02-07 21:48:59"
$ws.Range("H9").Value = "'74%"
$ws.Range("N9").Value = "4.1 °C 21:11 TU"
$ws.Range("O9").Value = "10.9 °C"
$ws.Range("E10").Value = "2026-02-07 21:49:02"
$ws.Range("O10").Value = "8.3 °C"
$ws.Range("E11").Value = "2026-02-07 21:49:05"
$ws.Range("E12").Value = "2026-02-07 21:49:07"
$ws.Range("O12").Value = "10.2 °C"
$ws.Range("E13").Value = "2026-02-07 21:49:09"
$ws.Range("E14").Value = "2026-02-07 21:49:11"
$ws.Range("H14").Value = "'60%"
$ws.Range("O14").Value = "11.8 °C"
$ws.Range("E15").Value = "2026-02-07 21:49:14"
$ws.Range("H15").Value = "'70%"
$ws.Range("N15").Value = "4.6 °C 21:00 TU"
$ws.Range("O15").Value = "10.5 °C"
$ws.Range("E16").Value = "2026-02-07 21:49:17"
$ws.Range("O16").Value = "-5.7 °C"
$ws.Range("E17").Value = "2026-02-07 21:49:19"
$ws.Range("L17").Value = "65.5 km/h - 242º 21:03 TU"
$ws.Range("O17").Value = "0.1 °C"
$ws.Range("E18").Value = "2026-02-07 21:49:22"
$ws.Range("J18").Value = "1004.3 hPa"
$ws.Range("E19").Value = "2026-02-07 21:49:24"
$ws.Range("E20").Value = "2026-02-07 21:49:27"
$ws.Range("E21").Value = "2026-02-07 21:49:30"
$ws.Range("H21").Value = "'83%"
$ws.Range("E22").Value = "2026-02-07 21:49:32"
$ws.Range("I22").Value = "1.5 mm"
$ws.Range("E23").Value = "2026-02-07 21:49:35"
$ws.Range("L23").Value = "44.3 km/h - 136º 21:10 TU"
$ws.Range("O23").Value = "-5.4 °C"
$ws.Range("E24").Value = "2026-02-07 21:49:38"
$ws.Range("I24").Value = "0.8 mm"
$ws.Range("J24").Value = "1007.2 hPa"
$ws.Range("E25").Value = "2026-02-07 21:49:40"
$ws.Range("H25").Value = "'79%"
$ws.Range("L25").Value = "22.0 km/h - 281º 21:24 TU"
$ws.Range("E26").Value = "2026-02-07 21:49:43"
$ws.Range("O26").Value = "2.8 °C"
$ws.Range("E27").Value = "2026-02-07 21:49:46"
$ws.Range("E28").Value = "2026-02-07 21:49:49"
$ws.Range("E29").Value = "2026-02-07 21:49:51"
$ws.Range("E30").Value = "2026-02-07 21:49:54"
$ws.Range("O30").Value = "9.8 °C"
$ws.Range("E31").Value = "2026-02-07 21:49:56"
$ws.Range("E32").Value = "2026-02-07 21:49:59"
$ws.Range("H32").Value = "'78%"
$ws.Range("I32").Value = "0.4 mm"
$ws.Range("E33").Value = "2026-02-07 21:50:01"
$ws.Range("J33").Value = "1006.2 hPa"
$ws.Range("N33").Value = "0.0 °C 21:22 TU"
$ws.Range("O33").Value = "2.2 °C"
$ws.Range("E34").Value = "2026-02-07 21:50:04"
$ws.Range("H34").Value = "'69%"
$ws.Range("L34").Value = "39.2 km/h - 286º 21:20 TU"
$ws.Range("E35").Value = "2026-02-07 21:50:07"
$ws.Range("J35").Value = "1007.4 hPa"
$ws.Range("L35").Value = "62.3 km/h - 213º 21:13 TU"
$ws.Range("E36").Value = "2026-02-07 21:50:09"
$ws.Range("K36").Value = "11.1 MJ/m2"
$ws.Range("O36").Value = "11.4 °C"
$ws.Range("E37").Value = "2026-02-07 21:50:12"
$ws.Range("O37").Value = "4.8 °C"
$ws.Range("E38").Value = "2026-02-07 21:50:15"
$ws.Range("O38").Value = "12.0 °C"
$ws.Range("E39").Value = "2026-02-07 21:50:17"
$ws.Range("L39").Value = "68.8 km/h - 303º 21:24 TU"
$ws.Range("E40").Value = "2026-02-07 21:50:20"
$ws.Range("E41").Value = "2026-02-07 21:50:22"
$ws.Range("H41").Value = "'56%"
$ws.Range("E42").Value = "2026-02-07 21:50:25"
$ws.Range("O42").Value = "10.3 °C"
$ws.Range("E43").Value = "2026-02-07 21:50:27"
$ws.Range("E44").Value = "2026-02-07 21:50:30"
$ws.Range("E45").Value = "2026-02-07 21:50:33"
$ws.Range("E46").Value = "2026-02-07 21:50:35"
$ws.Range("J46").Value = "1007.5 hPa"
$ws.Range("O46").Value = "9.1 °C"
